$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "s"
$ws.Range("B3").Value = "ss"
$ws.Range("C3").Value = "ss"

$ws.Range("G5").Select()
